$d = $word.ActiveDocument

$d.Content.Find.Execute("139 ± 258", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "139 ± 258 (74)", 2)

$d.Content.Find.Execute("93 ± 58", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "93 ± 58 (82)", 2)
